# TC09_C3DC_phs002518_LastKnSurStatus-Dead.xlsx
# "Updated remaining queries for C3DC"
#
# The C3DC SQL queries stored in column B/C of Sheet1 joined df_participant /
# df_diagnoses / df_treatments / df_treatment_resp / df_survival /
# df_reference_files using the generic `id` column. Those tables were
# renamed to use explicit `study_id` / `participant_id` columns, so every
# query's LEFT JOIN block needs to be updated to match. All affected cells
# share the exact same JOIN block, so we can apply one regex replace per
# cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldJoinBlock = @'
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_diagnoses dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_treatments trt ON prt.id = trt."participant.id"
LEFT JOIN 
    df_treatment_resp trr ON prt.id = trr."participant.id"
LEFT JOIN 
    df_survival srv ON prt.id = srv."participant.id"
LEFT JOIN 
    df_reference_files rfs ON std.id = rfs."study.id"
'@

$newJoinBlock = @'
LEFT JOIN 
    df_participant prt ON std.study_id = prt."study.study_id"
LEFT JOIN 
    df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"
LEFT JOIN 
    df_treatments trt ON prt.participant_id = trt."participant.participant_id"
LEFT JOIN 
    df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"
LEFT JOIN 
    df_survival srv ON prt.participant_id = srv."participant.participant_id"
LEFT JOIN 
    df_reference_files rfs ON std.study_id = rfs."study.study_id"
'@

$cellsWithQueries = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")

foreach ($addr in $cellsWithQueries) {
    $cell = $ws.Range($addr)
    $text = $cell.Text
    $updated = $text.Replace($oldJoinBlock, $newJoinBlock)
    $cell.Value = $updated
}

# Widen column C (the StatQuery column) from its old auto-fit width to a
# fixed 67.5 characters and drop the "best fit" auto-sizing flag. Excel's
# ColumnWidth setter re-derives the stored `width` from a pixel width that
# always differs from the requested character width by 5/6 of a character,
# so back that constant out to land exactly on 67.5 in the saved file.
$ws.Columns.Item(3).ColumnWidth = 67.5 - (5 / 6)

# Scroll the sheet view down one row so the visible viewport's top-left
# cell is A2 (selection stays on B2).
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
